# 3D CAD Stage 1 - Done up to Arm 1
# Updates BOM sheet: swap out old supplier links for DFRobot product pages,
# style the ball-caster link as a hyperlink, and move the stepper-motor URL
# into the correct column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 8 - Ball Casters: swap the pimoroni link for the dfrobot product page
# and format it like a hyperlink.
$ws.Range("F8").Value = "https://www.dfrobot.com/product-509.html"
$ws.Range("F8").Style = "Hyperlink"

# Row 4 - Stepper motors FIT0503: the reference URL was in G4 (digikey),
# move the (new) URL into F4 and clear the stray G4 cell.
$ws.Range("G4").ClearContents()
$ws.Range("F4").Value = "https://www.dfrobot.com/product-1508.html"

# Update the active selection on the BOM sheet to K11.
$ws.Range("K11").Select()
